$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" suffixed header columns to the respective
# input-file-name suffixes "_FV2410" / "_FV2504".
$headersLeft = @(
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headersLeft.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersLeft[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

$headersRight = @(
  "Segmentname_FV2504",
  "Segmentgruppe_FV2504",
  "Segment_FV2504",
  "Datenelement_FV2504",
  "Segment ID_FV2504",
  "Code_FV2504",
  "Qualifier_FV2504",
  "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504",
  "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headersRight.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersRight[$i]
}

# Turn the data range into an Excel Table (ListObject) with headers taken
# from row 1 (which now carries the renamed column titles above).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U74"), 0, 1)
$tbl.Name = "Table1"

# Freeze the header row (row 1) so it stays visible while scrolling.
$null = $ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
